# Ultrasonic Sensor pin change
# The Trigger/Echo rows move from GP19/GP18 (rows 37-38) down to GP17/GP16
# (rows 40-41); the GP19/GP18 and GND rows keep only their Pin (col B) value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37 (GP19) and 38 (GP18): drop the Name/Description, keep Pin.
$ws.Range("A37").ClearContents()
$ws.Range("C37").ClearContents()
$ws.Range("A38").ClearContents()
$ws.Range("C38").ClearContents()

# Row 40 (GP17) and 41 (GP16): now hold the Ultrasonic Sensor Trigger/Echo
# Name + Description, keeping their existing Pin value.
$ws.Range("A40").Value = "Ultrasonic Sensor Trigger"
$ws.Range("C40").Value = "Sends the ultrasonic pulse"
$ws.Range("A41").Value = "Ultrasonic Sensor Echo"
$ws.Range("C41").Value = "Recieves the returning ultrasonic signal"

# View-state updates to match the author's saved selection/scroll position.
$ws.Columns.Item(1).ColumnWidth = 23

$ws.Activate()
$ws.Range("C38").Select()

$win = $excel.ActiveWindow
$win.ScrollRow = 29
$win.ScrollColumn = 1
